# CMSIS Core (Using TrustZone) documentation update.
#
# Fixes the label on the "Snip Single Corner Rectangle 23" shape (inside
# "Group 20") on slide 5 from "partitions_<device>.h" to the correct
# header-file name "partition_<device>.h" (the project file is actually
# named partition_<device>.h, not partitions_<device>.h).

$msoGroup = 6

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
        if ($shp.Type -eq $msoGroup) {
            $found = Find-ShapeByName $shp.GroupItems $name
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$target = Find-ShapeByName $s.Shapes "Snip Single Corner Rectangle 23"

$tr = $target.TextFrame.TextRange

# The paragraph currently holds two runs that together spell
# "p" + "artitions_<device>.h" = "partitions_<device>.h".
# Re-split the same characters into "partition_<" + "device>.h", updating
# the second (longer) run first so the first run's character offset
# doesn't shift underneath us.
$run2 = $tr.Characters(2, 21)
$run2.Text = "device>.h"

$run1 = $tr.Characters(1, 1)
$run1.Text = "partition_<"
